$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - update values only
$ws.Range("B3").Value = 0.9991150998960281
$ws.Range("C3").Value = 0.9991567441848135
$ws.Range("D3").Value = 0.996933037205515

# Row 4: label change GradientBoostingRegressor -> DecisionTreeRegressor, update values
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9990874661804025
$ws.Range("C4").Value = 0.9991447349734291
$ws.Range("D4").Value = 0.9981701652541393

# Row 5: label change AdaBoostRegressor -> MLPRegressor, update values
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.999283417342841
$ws.Range("C5").Value = 0.9993132564651125
$ws.Range("D5").Value = 0.999341745970114
